$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "RM 232" row (row 26) and the "SC 92" row (originally row 28,
# which becomes row 27 once row 26 is removed). Deleting shifts everything
# below up, matching the diff (dimension goes from A1:F35 to A1:F33).
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# Toggle individual F-column (error) values: some get filled in, some get
# cleared out, per the diff.
$ws.Range("F3").Value = 17.64
$ws.Range("F5").Value = ""
$ws.Range("F21").Value = 16.58
$ws.Range("F23").Value = ""

# "SC 193" (now row 32 after the two row deletions) previously had a blank
# F value; the diff shows it filled in with 17.39.
$ws.Range("F32").Value = 17.39
